# Update cryptocurrency price/volume figures per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextValue "D2" "41.996.86"
Set-TextValue "E2" "  -0.59%  "
Set-TextValue "D3" "2.213.68"
Set-TextValue "E3" "  -1.50%  "
Set-TextValue "E4" "  +0.23%  "
Set-TextValue "D5" "241.48"
Set-TextValue "E5" "  -2.34%  "
Set-TextValue "D6" "0.626"
Set-TextValue "E6" "  +0.74%  "
Set-TextValue "D7" "73.39"
Set-TextValue "E7" "  -1.69%  "
Set-TextValue "E8" "  +0.09%  "
Set-TextValue "D9" "0.609"
Set-TextValue "E9" "  -1.32%  "
Set-TextValue "D10" "43.41"
Set-TextValue "E10" "  +2.42%  "
Set-TextValue "D11" "0.0956"
Set-TextValue "E11" "  +1.66%  "
Set-TextValue "E12" "  -1.51%  "
Set-TextValue "E13" "  -0.28%  "
Set-TextValue "D14" "2.546.39"
Set-TextValue "E14" "  -1.39%  "
Set-TextValue "D15" "14.22"
Set-TextValue "E15" "  -2.38%  "
Set-TextValue "E16" "  -1.30%  "
Set-TextValue "D17" "2.220.16"
Set-TextValue "E17" "  -1.61%  "
Set-TextValue "D18" "41.810.75"
Set-TextValue "E18" "  -0.76%  "
Set-TextValue "D19" "0.0000109"
Set-TextValue "E19" "  +10.20%  "
Set-TextValue "D20" "72.47"
Set-TextValue "E20" "  +0.67%  "
Set-TextValue "D21" "6.12"
Set-TextValue "E21" "  -0.34%  "
Set-TextValue "D22" "10.37"
Set-TextValue "E22" "  +15.88%  "
Set-TextValue "D23" "228.92"
Set-TextValue "E23" "  -0.88%  "
Set-TextValue "E24" "  -6.70%  "
Set-TextValue "E25" "  +0.05%  "
Set-TextValue "D26" "11.46"
Set-TextValue "E26" "  +1.52%  "
Set-TextValue "E27" "  -0.71%  "
Set-TextValue "D28" "2.27"
Set-TextValue "E28" "  -1.91%  "
Set-TextValue "E29" "  -0.40%  "
Set-TextValue "D30" "166.90"
Set-TextValue "E30" "  -1.43%  "
Set-TextValue "D31" "20.55"
Set-TextValue "E31" "  -0.71%  "
Set-TextValue "D32" "5.53"
Set-TextValue "E32" "  +4.98%  "
Set-TextValue "E33" "  -4.30%  "
Set-TextValue "E34" "  +0.43%  "
Set-TextValue "D35" "28.91"
Set-TextValue "E35" "  -4.88%  "
Set-TextValue "D36" "0.111"
Set-TextValue "E36" "  -7.26%  "
Set-TextValue "E37" "  -6.11%  "
Set-TextValue "E38" "  -1.03%  "
Set-TextValue "D39" "12.84"
Set-TextValue "E39" "  -4.76%  "
Set-TextValue "D40" "65.81"
Set-TextValue "E40" "  +6.26%  "
Set-TextValue "E41" "  -3.45%  "
Set-TextValue "E42" "  -3.47%  "
Set-TextValue "D43" "0.199"
Set-TextValue "E43" "  -1.53%  "
Set-TextValue "D44" "8.67"
Set-TextValue "E44" "  -0.05%  "
Set-TextValue "D45" "103.67"
Set-TextValue "E45" "  -4.83%  "
Set-TextValue "E46" "  -1.14%  "
Set-TextValue "D47" "2.42"
Set-TextValue "E47" "  +4.81%  "
Set-TextValue "E48" "  -1.21%  "
Set-TextValue "E49" "  -0.55%  "
Set-TextValue "D50" "2.70"
Set-TextValue "E50" "  +0.26%  "
Set-TextValue "D51" "2.419.34"
Set-TextValue "E51" "  -1.48%  "
